$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-match data (everything except id/Div/Div Original
# Name/Date, i.e. B and F through AC).
$dataCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $dataCols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# ---------------------------------------------------------------------------
# Rows 78 <-> 79 : full data swap (column A keeps its own row index).
# ---------------------------------------------------------------------------
Swap-Rows 78 79

# ---------------------------------------------------------------------------
# Rows 87 <-> 88 : full data swap.
# ---------------------------------------------------------------------------
Swap-Rows 87 88

# ---------------------------------------------------------------------------
# Rows 237 <-> 238 : full data swap.
# ---------------------------------------------------------------------------
Swap-Rows 237 238

# ---------------------------------------------------------------------------
# Rows 264-266 : refreshed match data for the league update, rows 267 and
# 268 are removed outright (the raw feed only has 3 new fixtures here).
# ---------------------------------------------------------------------------
$ws.Range("B264").Value = 6992669
$ws.Range("E264").Value = 45347.3125
$ws.Range("F264").Value = "Ratchaburi FC"
$ws.Range("G264").Value = "Uthai Thani FC"
$ws.Range("K264").Value = 1.75
$ws.Range("L264").Value = 3.5
$ws.Range("M264").Value = 4
$ws.Range("N264").Value = 1.7
$ws.Range("O264").Value = 3.6
$ws.Range("P264").Value = 4.2
$ws.Range("Q264").Value = -0.75
$ws.Range("R264").Value = 1.9
$ws.Range("S264").Value = 1.9
$ws.Range("T264").Value = 2.75
$ws.Range("U264").Value = 1.975
$ws.Range("V264").Value = 1.825

$ws.Range("B265").Value = 6992673
$ws.Range("E265").Value = 45347.33333333334
$ws.Range("F265").Value = "BG Pathum United"
$ws.Range("G265").Value = "Bangkok United"
$ws.Range("K265").Value = 2.8
$ws.Range("L265").Value = 3.4
$ws.Range("M265").Value = 2.2
$ws.Range("N265").Value = 2.375
$ws.Range("O265").Value = 3.3
$ws.Range("P265").Value = 2.625
$ws.Range("Q265").Value = 0
$ws.Range("R265").Value = 1.825
$ws.Range("S265").Value = 1.975
$ws.Range("T265").Value = 2.5
$ws.Range("U265").Value = 1.875
$ws.Range("V265").Value = 1.925

$ws.Range("B266").Value = 6992670
$ws.Range("E266").Value = 45347.375
$ws.Range("F266").Value = "Chonburi"
$ws.Range("G266").Value = "Sukhothai FC"
$ws.Range("K266").Value = 1.95
$ws.Range("L266").Value = 3.5
$ws.Range("M266").Value = 3.25
$ws.Range("N266").Value = 1.8
$ws.Range("O266").Value = 3.6
$ws.Range("P266").Value = 3.6
$ws.Range("Q266").Value = -0.5
$ws.Range("R266").Value = 1.8
$ws.Range("S266").Value = 2
$ws.Range("T266").Value = 2.75
$ws.Range("U266").Value = 1.9
$ws.Range("V266").Value = 1.9

$ws.Rows("267:268").Delete()
